$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2825.3103
$ws.Range("J17").Value = 2400.2964
$ws.Range("L17").Value = 7200.889200000001
$ws.Range("N17").Value = -7536.889200000001
$ws.Range("H64").Value = 3753.75
$ws.Range("J64").Value = 4750
$ws.Range("L64").Value = 4750
$ws.Range("N64").Value = -5246
$ws.Range("H67").Value = 3753.75
$ws.Range("J67").Value = 4750
$ws.Range("L67").Value = 4750
$ws.Range("N67").Value = -6466
$ws.Range("H88").Value = 35715984
$ws.Range("I88").Value = 125001330
$ws.Range("J88").Value = 1849.6
$ws.Range("K88").Value = 125001330
$ws.Range("L88").Value = 1849.6
$ws.Range("M88").Value = -125000924
$ws.Range("N88").Value = -2661.6
$ws.Range("H91").Value = 35715984
$ws.Range("I91").Value = 125001330
$ws.Range("J91").Value = 1849.6
$ws.Range("K91").Value = 125001330
$ws.Range("L91").Value = 1849.6
$ws.Range("M91").Value = -124999926
$ws.Range("N91").Value = -4657.6
$ws.Range("H129").Value = 865.89746
$ws.Range("I129").Value = 638.6
$ws.Range("K129").Value = 1915.8
$ws.Range("M129").Value = 3084.2
$ws.Range("H138").Value = 1875.2715
$ws.Range("I138").Value = 1702.5588
$ws.Range("K138").Value = 5107.6764
$ws.Range("M138").Value = 32.32359999999971

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 695739.6
$ws.Range("I2").Value = 926984
$ws.Range("K2").Value = 926984
$ws.Range("M2").Value = -926871
$ws.Range("H37").Value = 16264
$ws.Range("I37").Value = 20048
$ws.Range("J37").Value = 12480
$ws.Range("K37").Value = 20048
$ws.Range("L37").Value = 12480
$ws.Range("M37").Value = -19775
$ws.Range("N37").Value = -13026
$ws.Range("H110").Value = 141.16667
$ws.Range("I110").Value = 141.16667
$ws.Range("K110").Value = 141.16667
$ws.Range("M110").Value = 1903.83333
$ws.Range("H116").Value = 695739.6
$ws.Range("I116").Value = 926984
$ws.Range("K116").Value = 926984
$ws.Range("M116").Value = -924690

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 695739.6
$ws.Range("I3").Value = 926984
$ws.Range("K3").Value = 926984
$ws.Range("M3").Value = -926870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2261
$ws.Range("I86").Value = 2413.2
$ws.Range("K86").Value = 2413.2
$ws.Range("M86").Value = -1290.2
$ws.Range("H89").Value = 2261
$ws.Range("I89").Value = 2413.2
$ws.Range("K89").Value = 12066
$ws.Range("M89").Value = -6450
$ws.Range("H94").Value = 1088.6364
$ws.Range("I94").Value = 982
$ws.Range("K94").Value = 982
$ws.Range("M94").Value = -531
$ws.Range("H99").Value = 3091.1
$ws.Range("I99").Value = 2164.6
$ws.Range("K99").Value = 2164.6
$ws.Range("M99").Value = -666.5999999999999
$ws.Range("H107").Value = 754.7059
$ws.Range("I107").Value = 533.4286
$ws.Range("J107").Value = 1787.3334
$ws.Range("K107").Value = 533.4286
$ws.Range("L107").Value = 1787.3334
$ws.Range("M107").Value = 1386.5714
$ws.Range("N107").Value = -5627.3334
$ws.Range("H126").Value = 3091.1
$ws.Range("I126").Value = 2164.6
$ws.Range("K126").Value = 6493.799999999999
$ws.Range("M126").Value = -4023.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 203.36363
$ws.Range("J2").Value = 265.25
$ws.Range("L2").Value = 1591.5
$ws.Range("N2").Value = -1817.5
$ws.Range("H10").Value = 256.5
$ws.Range("I10").Value = 256.5
$ws.Range("K10").Value = 769.5
$ws.Range("M10").Value = -630.5
$ws.Range("H109").Value = 4993.4
$ws.Range("I109").Value = 1309.6666
$ws.Range("J109").Value = 5914.3335
$ws.Range("K109").Value = 3928.9998
$ws.Range("L109").Value = 17743.0005
$ws.Range("M109").Value = -2888.9998
$ws.Range("N109").Value = -19823.0005
$ws.Range("H134").Value = 2697.5
$ws.Range("I134").Value = 1597.0588
$ws.Range("J134").Value = 4398.1816
$ws.Range("K134").Value = 4791.1764
$ws.Range("L134").Value = 13194.5448
$ws.Range("M134").Value = 278.8235999999997
$ws.Range("N134").Value = -23334.5448
$ws.Range("H140").Value = 3291.5557
$ws.Range("I140").Value = 1627.8889
$ws.Range("J140").Value = 4955.222
$ws.Range("K140").Value = 4883.6667
$ws.Range("L140").Value = 14865.666
$ws.Range("M140").Value = 296.3333000000002
$ws.Range("N140").Value = -25225.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1166.7
$ws.Range("I113").Value = 1118
$ws.Range("J113").Value = 1187.5714
$ws.Range("K113").Value = 1118
$ws.Range("L113").Value = 1187.5714
$ws.Range("M113").Value = 1052
$ws.Range("N113").Value = -5527.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2109.3635
$ws.Range("I7").Value = 1817.0625
$ws.Range("K7").Value = 1817.0625
$ws.Range("M7").Value = -1705.0625
$ws.Range("H93").Value = 15874003
$ws.Range("I93").Value = 851.5
$ws.Range("K93").Value = 851.5
$ws.Range("M93").Value = 396.5
$ws.Range("H122").Value = 6344.222
$ws.Range("I122").Value = 6639.3
$ws.Range("K122").Value = 19917.9
$ws.Range("M122").Value = -17467.9
$ws.Range("H126").Value = 2109.3635
$ws.Range("I126").Value = 1817.0625
$ws.Range("K126").Value = 5451.1875
$ws.Range("M126").Value = -2981.1875
$ws.Range("H132").Value = 2597.9443
$ws.Range("I132").Value = 1899.5
$ws.Range("K132").Value = 5698.5
$ws.Range("M132").Value = -3168.5
$ws.Range("H136").Value = 3887.9375
$ws.Range("I136").Value = 3153.08
$ws.Range("J136").Value = 6512.4287
$ws.Range("K136").Value = 9459.24
$ws.Range("L136").Value = 19537.2861
$ws.Range("M136").Value = -6909.24
$ws.Range("N136").Value = -24637.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1571.375
$ws.Range("I96").Value = 990
$ws.Range("J96").Value = 1654.4286
$ws.Range("K96").Value = 990
$ws.Range("L96").Value = 1654.4286
$ws.Range("N96").Value = -4400.4286
$ws.Range("M96").Value = 383
$ws.Range("H100").Value = 1076.3334
$ws.Range("I100").Value = 931.6
$ws.Range("K100").Value = 1863.2
$ws.Range("M100").Value = -1322.2
$ws.Range("H122").Value = 61358.77
$ws.Range("I122").Value = 98596.125
$ws.Range("K122").Value = 295788.375
$ws.Range("M122").Value = -293338.375
$ws.Range("H126").Value = 2174.3333
$ws.Range("I126").Value = 1796.6666
$ws.Range("J126").Value = 4062.6667
$ws.Range("K126").Value = 5389.9998
$ws.Range("L126").Value = 12188.0001
$ws.Range("M126").Value = -2919.9998
$ws.Range("N126").Value = -17128.0001
$ws.Range("H136").Value = 16341860
$ws.Range("I136").Value = 27779786
$ws.Range("J136").Value = 1964.1428
$ws.Range("K136").Value = 83339358
$ws.Range("L136").Value = 5892.428400000001
$ws.Range("M136").Value = -83336808
$ws.Range("N136").Value = -10992.4284
